# Apply the "add 2022-Q4 data" edit:
#  1. Insert a brand-new worksheet named "2022-Q4" right after "总计" (i.e.
#     right before "2022-Q3"), and populate it with the two new fund rows.
#  2. Insert a new row into the "总计" (summary) sheet right after its
#     header row, holding the 2022-Q4 aggregate (count=2, value=0.93), and
#     renumber the leading index column (A) for the rows that shifted down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q4" worksheet
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row (bold, centered, bordered - matches the other quarter sheets)
$hdr = $q4.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Leading index column (A2:A3) - bold, centered, bordered
$idx = $q4.Range("A2:A3")
$idx.Font.Bold = $true
$idx.HorizontalAlignment = -4108
$idx.VerticalAlignment = -4160
$idx.Borders.LineStyle = 1
$idx.Borders.Weight = 2

$q4.Range("A2").Value = 0
$q4.Range("A3").Value = 1

# Data rows - force text storage ("@") so fund codes keep leading zeros and
# numeric-looking strings (percentages, AUM figures) aren't coerced to
# binary doubles.
$q4.Range("B2:G3").NumberFormat = "@"

$q4.Range("B2").Value = "002121"
$q4.Range("C2").Value = "广发沪港深新起点股票A"
$q4.Range("D2").Value = "26.30"
$q4.Range("E2").Value = "88.97"
$q4.Range("F2").Value = "3.46"
$q4.Range("G2").Value = "0.9100"
$q4.Range("H2").Value = 9

$q4.Range("B3").Value = "010024"
$q4.Range("C3").Value = "广发沪港深新起点股票C"
$q4.Range("D3").Value = "0.59"
$q4.Range("E3").Value = "88.97"
$q4.Range("F3").Value = "3.46"
$q4.Range("G3").Value = "0.0204"
$q4.Range("H3").Value = 9

# ---------------------------------------------------------------------
# 2. "总计" (summary) sheet - insert the 2022-Q4 row
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows(2).Insert()

$a2 = $summary.Range("A2")
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1
$a2.Borders.Weight = 2
$a2.Value = 0

$summary.Range("B2").NumberFormat = "@"
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.93

# Renumber the index column for the rows that shifted down one position
# (they keep their original quarter data, just a new running index).
for ($i = 0; $i -le 6; $i++) {
    $row = 3 + $i
    $summary.Range("A$row").Value = $i + 1
}
